$d = $word.ActiveDocument

$d.Content.Find.Execute('2024-12-31 Tuesday', $true, $false, $false, $false, $false, $true, 1, $false, '2025-01-01 Wednesday', 2) | Out-Null
$d.Content.Find.Execute('76-29=47', $true, $false, $false, $false, $false, $true, 1, $false, '0+18=18', 2) | Out-Null
$d.Content.Find.Execute('7-3=4', $true, $false, $false, $false, $false, $true, 1, $false, '56+29=85', 2) | Out-Null
$d.Content.Find.Execute('15+32=47', $true, $false, $false, $false, $false, $true, 1, $false, '84-64=20', 2) | Out-Null
$d.Content.Find.Execute('27+17=44', $true, $false, $false, $false, $false, $true, 1, $false, '4+15=19', 2) | Out-Null
$d.Content.Find.Execute('99-6=93', $true, $false, $false, $false, $false, $true, 1, $false, '14+46=60', 2) | Out-Null
$d.Content.Find.Execute('50-7=43', $true, $false, $false, $false, $false, $true, 1, $false, '78-35=43', 2) | Out-Null
$d.Content.Find.Execute('55-32=23', $true, $false, $false, $false, $false, $true, 1, $false, '2+79=81', 2) | Out-Null
$d.Content.Find.Execute('44-15=29', $true, $false, $false, $false, $false, $true, 1, $false, '48-17=31', 2) | Out-Null
$d.Content.Find.Execute('81-11=70', $true, $false, $false, $false, $false, $true, 1, $false, '30+9=39', 2) | Out-Null
$d.Content.Find.Execute('57-46=11', $true, $false, $false, $false, $false, $true, 1, $false, '0+5=5', 2) | Out-Null
$d.Content.Find.Execute('47+40=87', $true, $false, $false, $false, $false, $true, 1, $false, '28+37=65', 2) | Out-Null
$d.Content.Find.Execute('13+73=86', $true, $false, $false, $false, $false, $true, 1, $false, '82-47=35', 2) | Out-Null
$d.Content.Find.Execute('80-76=4', $true, $false, $false, $false, $false, $true, 1, $false, '62+31=93', 2) | Out-Null
$d.Content.Find.Execute('0+22=22', $true, $false, $false, $false, $false, $true, 1, $false, '86-67=19', 2) | Out-Null
$d.Content.Find.Execute('26+42=68', $true, $false, $false, $false, $false, $true, 1, $false, '86-73=13', 2) | Out-Null
$d.Content.Find.Execute('57-49=8', $true, $false, $false, $false, $false, $true, 1, $false, '22+22=44', 2) | Out-Null
$d.Content.Find.Execute('69-18=51', $true, $false, $false, $false, $false, $true, 1, $false, '29+61=90', 2) | Out-Null
$d.Content.Find.Execute('54-15=39', $true, $false, $false, $false, $false, $true, 1, $false, '90-64=26', 2) | Out-Null
$d.Content.Find.Execute('50-31=19', $true, $false, $false, $false, $false, $true, 1, $false, '68-37=31', 2) | Out-Null
$d.Content.Find.Execute('75-36=39', $true, $false, $false, $false, $false, $true, 1, $false, '85-23=62', 2) | Out-Null
$d.Content.Find.Execute('93-40=53', $true, $false, $false, $false, $false, $true, 1, $false, '96-64=32', 2) | Out-Null
$d.Content.Find.Execute('43+47=90', $true, $false, $false, $false, $false, $true, 1, $false, '56+11=67', 2) | Out-Null
$d.Content.Find.Execute('46+45=91', $true, $false, $false, $false, $false, $true, 1, $false, '98-84=14', 2) | Out-Null
$d.Content.Find.Execute('60-46=14', $true, $false, $false, $false, $false, $true, 1, $false, '48-26=22', 2) | Out-Null
$d.Content.Find.Execute('23+65=88', $true, $false, $false, $false, $false, $true, 1, $false, '6+14=20', 2) | Out-Null
$d.Content.Find.Execute('65-1=64', $true, $false, $false, $false, $false, $true, 1, $false, '13+80=93', 2) | Out-Null
$d.Content.Find.Execute('88-56=32', $true, $false, $false, $false, $false, $true, 1, $false, '2+66=68', 2) | Out-Null
$d.Content.Find.Execute('8+11=19', $true, $false, $false, $false, $false, $true, 1, $false, '37+61=98', 2) | Out-Null
$d.Content.Find.Execute('2+7=9', $true, $false, $false, $false, $false, $true, 1, $false, '60-29=31', 2) | Out-Null
$d.Content.Find.Execute('22+58=80', $true, $false, $false, $false, $false, $true, 1, $false, '32+33=65', 2) | Out-Null
$d.Content.Find.Execute('3+24=27', $true, $false, $false, $false, $false, $true, 1, $false, '39+51=90', 2) | Out-Null
$d.Content.Find.Execute('67-28=39', $true, $false, $false, $false, $false, $true, 1, $false, '78-41=37', 2) | Out-Null
$d.Content.Find.Execute('11+68=79', $true, $false, $false, $false, $false, $true, 1, $false, '46-7=39', 2) | Out-Null
$d.Content.Find.Execute('57-26=31', $true, $false, $false, $false, $false, $true, 1, $false, '23-15=8', 2) | Out-Null
$d.Content.Find.Execute('61+13=74', $true, $false, $false, $false, $false, $true, 1, $false, '97-4=93', 2) | Out-Null
$d.Content.Find.Execute('42+15=57', $true, $false, $false, $false, $false, $true, 1, $false, '35-3=32', 2) | Out-Null
$d.Content.Find.Execute('76+2=78', $true, $false, $false, $false, $false, $true, 1, $false, '88-11=77', 2) | Out-Null
$d.Content.Find.Execute('38-5=33', $true, $false, $false, $false, $false, $true, 1, $false, '17-10=7', 2) | Out-Null
$d.Content.Find.Execute('8+36=44', $true, $false, $false, $false, $false, $true, 1, $false, '41-10=31', 2) | Out-Null
$d.Content.Find.Execute('2+48=50', $true, $false, $false, $false, $false, $true, 1, $false, '6+16=22', 2) | Out-Null
$d.Content.Find.Execute('5+52=57', $true, $false, $false, $false, $false, $true, 1, $false, '57-40=17', 2) | Out-Null
$d.Content.Find.Execute('47+46=93', $true, $false, $false, $false, $false, $true, 1, $false, '55+4=59', 2) | Out-Null
$d.Content.Find.Execute('97-22=75', $true, $false, $false, $false, $false, $true, 1, $false, '28+4=32', 2) | Out-Null
$d.Content.Find.Execute('10-4=6', $true, $false, $false, $false, $false, $true, 1, $false, '91+6=97', 2) | Out-Null
$d.Content.Find.Execute('32-26=6', $true, $false, $false, $false, $false, $true, 1, $false, '99-59=40', 2) | Out-Null
$d.Content.Find.Execute('15-8=7', $true, $false, $false, $false, $false, $true, 1, $false, '87-5=82', 2) | Out-Null
$d.Content.Find.Execute('86-61=25', $true, $false, $false, $false, $false, $true, 1, $false, '23+72=95', 2) | Out-Null
$d.Content.Find.Execute('24+74=98', $true, $false, $false, $false, $false, $true, 1, $false, '98-16=82', 2) | Out-Null
$d.Content.Find.Execute('89+10=99', $true, $false, $false, $false, $false, $true, 1, $false, '87+3=90', 2) | Out-Null
$d.Content.Find.Execute('0+40=40', $true, $false, $false, $false, $false, $true, 1, $false, '8+1=9', 2) | Out-Null
$d.Content.Find.Execute('30-27=3', $true, $false, $false, $false, $false, $true, 1, $false, '91-47=44', 2) | Out-Null
$d.Content.Find.Execute('2+5=7', $true, $false, $false, $false, $false, $true, 1, $false, '93-34=59', 2) | Out-Null
$d.Content.Find.Execute('91-60=31', $true, $false, $false, $false, $false, $true, 1, $false, '75+22=97', 2) | Out-Null
$d.Content.Find.Execute('61-38=23', $true, $false, $false, $false, $false, $true, 1, $false, '6+15=21', 2) | Out-Null
$d.Content.Find.Execute('96+0=96', $true, $false, $false, $false, $false, $true, 1, $false, '69-3=66', 2) | Out-Null
$d.Content.Find.Execute('78-65=13', $true, $false, $false, $false, $false, $true, 1, $false, '48+40=88', 2) | Out-Null
$d.Content.Find.Execute('89+0=89', $true, $false, $false, $false, $false, $true, 1, $false, '89-42=47', 2) | Out-Null
$d.Content.Find.Execute('61-34=27', $true, $false, $false, $false, $false, $true, 1, $false, '8+6=14', 2) | Out-Null
$d.Content.Find.Execute('45+5=50', $true, $false, $false, $false, $false, $true, 1, $false, '31+13=44', 2) | Out-Null
$d.Content.Find.Execute('54-36=18', $true, $false, $false, $false, $false, $true, 1, $false, '97-88=9', 2) | Out-Null
$d.Content.Find.Execute('23-3=20', $true, $false, $false, $false, $false, $true, 1, $false, '97-33=64', 2) | Out-Null
$d.Content.Find.Execute('92-86=6', $true, $false, $false, $false, $false, $true, 1, $false, '55-21=34', 2) | Out-Null
$d.Content.Find.Execute('85-35=50', $true, $false, $false, $false, $false, $true, 1, $false, '69-65=4', 2) | Out-Null
$d.Content.Find.Execute('35+35=70', $true, $false, $false, $false, $false, $true, 1, $false, '64-13=51', 2) | Out-Null
$d.Content.Find.Execute('76-8=68', $true, $false, $false, $false, $false, $true, 1, $false, '69-51=18', 2) | Out-Null
$d.Content.Find.Execute('25-21=4', $true, $false, $false, $false, $false, $true, 1, $false, '95-75=20', 2) | Out-Null
$d.Content.Find.Execute('37-37=0', $true, $false, $false, $false, $false, $true, 1, $false, '73+19=92', 2) | Out-Null
$d.Content.Find.Execute('60+6=66', $true, $false, $false, $false, $false, $true, 1, $false, '86-16=70', 2) | Out-Null
$d.Content.Find.Execute('99-36=63', $true, $false, $false, $false, $false, $true, 1, $false, '14+0=14', 2) | Out-Null
$d.Content.Find.Execute('34+21=55', $true, $false, $false, $false, $false, $true, 1, $false, '38+25=63', 2) | Out-Null
$d.Content.Find.Execute('16+30=46', $true, $false, $false, $false, $false, $true, 1, $false, '73-6=67', 2) | Out-Null
$d.Content.Find.Execute('63-26=37', $true, $false, $false, $false, $false, $true, 1, $false, '52-42=10', 2) | Out-Null
$d.Content.Find.Execute('53-41=12', $true, $false, $false, $false, $false, $true, 1, $false, '8+32=40', 2) | Out-Null
$d.Content.Find.Execute('35-31=4', $true, $false, $false, $false, $false, $true, 1, $false, '7+0=7', 2) | Out-Null
$d.Content.Find.Execute('29+55=84', $true, $false, $false, $false, $false, $true, 1, $false, '42+53=95', 2) | Out-Null
$d.Content.Find.Execute('6+87=93', $true, $false, $false, $false, $false, $true, 1, $false, '0+62=62', 2) | Out-Null
$d.Content.Find.Execute('40+47=87', $true, $false, $false, $false, $false, $true, 1, $false, '90-52=38', 2) | Out-Null
$d.Content.Find.Execute('78-20=58', $true, $false, $false, $false, $false, $true, 1, $false, '62-35=27', 2) | Out-Null
$d.Content.Find.Execute('12+61=73', $true, $false, $false, $false, $false, $true, 1, $false, '33+1=34', 2) | Out-Null
$d.Content.Find.Execute('37+59=96', $true, $false, $false, $false, $false, $true, 1, $false, '77-40=37', 2) | Out-Null
$d.Content.Find.Execute('13+26=39', $true, $false, $false, $false, $false, $true, 1, $false, '94-1=93', 2) | Out-Null
$d.Content.Find.Execute('19-11=8', $true, $false, $false, $false, $false, $true, 1, $false, '15+81=96', 2) | Out-Null
$d.Content.Find.Execute('47+13=60', $true, $false, $false, $false, $false, $true, 1, $false, '16+21=37', 2) | Out-Null
$d.Content.Find.Execute('98-53=45', $true, $false, $false, $false, $false, $true, 1, $false, '15+56=71', 2) | Out-Null
$d.Content.Find.Execute('44+21=65', $true, $false, $false, $false, $false, $true, 1, $false, '57+8=65', 2) | Out-Null
$d.Content.Find.Execute('56-34=22', $true, $false, $false, $false, $false, $true, 1, $false, '61-21=40', 2) | Out-Null
$d.Content.Find.Execute('29+58=87', $true, $false, $false, $false, $false, $true, 1, $false, '79-40=39', 2) | Out-Null
$d.Content.Find.Execute('95-27=68', $true, $false, $false, $false, $false, $true, 1, $false, '42+33=75', 2) | Out-Null
$d.Content.Find.Execute('82-76=6', $true, $false, $false, $false, $false, $true, 1, $false, '30+34=64', 2) | Out-Null
$d.Content.Find.Execute('93-84=9', $true, $false, $false, $false, $false, $true, 1, $false, '67+3=70', 2) | Out-Null
$d.Content.Find.Execute('94-92=2', $true, $false, $false, $false, $false, $true, 1, $false, '72+0=72', 2) | Out-Null
$d.Content.Find.Execute('53-45=8', $true, $false, $false, $false, $false, $true, 1, $false, '78-34=44', 2) | Out-Null
$d.Content.Find.Execute('37-30=7', $true, $false, $false, $false, $false, $true, 1, $false, '30+12=42', 2) | Out-Null
$d.Content.Find.Execute('90+2=92', $true, $false, $false, $false, $false, $true, 1, $false, '34+63=97', 2) | Out-Null
$d.Content.Find.Execute('70-30=40', $true, $false, $false, $false, $false, $true, 1, $false, '3+96=99', 2) | Out-Null
$d.Content.Find.Execute('26+66=92', $true, $false, $false, $false, $false, $true, 1, $false, '23+1=24', 2) | Out-Null
$d.Content.Find.Execute('33+52=85', $true, $false, $false, $false, $false, $true, 1, $false, '93-59=34', 2) | Out-Null
$d.Content.Find.Execute('14+4=18', $true, $false, $false, $false, $false, $true, 1, $false, '92-90=2', 2) | Out-Null
$d.Content.Find.Execute('91-49=42', $true, $false, $false, $false, $false, $true, 1, $false, '96-39=57', 2) | Out-Null
$d.Content.Find.Execute('71+23=94', $true, $false, $false, $false, $false, $true, 1, $false, '21+60=81', 2) | Out-Null
